# Update the "Förändrad" (Changed) date column (C) for data rows 2-185
# from 2023-09-09 (serial 45178) to 2023-09-10 (serial 45179).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 185; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45178) {
        $cell.Value2 = 45179
    }
}
